$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that currently sits in the body
#    right after the "Location: " run.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Split the single header/footer pair into even / default / first
#    variants (this is what Word does internally, minting header2.xml,
#    header3.xml, footer2.xml, footer3.xml, footer1... parts and wiring
#    up the extra <w:headerReference>/<w:footerReference> entries).
#    Touching the even & first Range objects is what forces Word to
#    actually materialise those parts instead of just lazily reporting
#    Exists = True.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)

$hdrs = $sec.Headers
$ftrs = $sec.Footers

$hdrFirst = $hdrs.Item(2)
$hdrEven  = $hdrs.Item(3)
$ftrFirst = $ftrs.Item(2)
$ftrEven  = $ftrs.Item(3)

$hdrFirst.Range.Text = ""
$hdrEven.Range.Text  = ""
$ftrFirst.Range.Text = ""
$ftrEven.Range.Text  = ""

# ------------------------------------------------------------------
# 3) In the primary/default header ("Buildit(TM)"), add a fresh
#    "_GoBack" bookmark right after the "Buildit" run (this has to
#    happen before the "TM" run is deleted so it lands between the
#    run and the spellEnd proof-error marker, matching the target),
#    then strip the superscript "TM" run entirely.
# ------------------------------------------------------------------
$hdrPrimary = $hdrs.Item(1)

$builditRange = $hdrPrimary.Range.Duplicate
$builditRange.Find.Execute("Buildit", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$builditRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $builditRange)

$tmRange = $hdrPrimary.Range.Duplicate
$tmRange.Find.Execute("TM", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "done"
